$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I86").Value = 66672230
$ws.Range("J86").Value = 6271.75
$ws.Range("K86").Value = 66672230
$ws.Range("L86").Value = 6271.75
$ws.Range("M86").Value = -66671107
$ws.Range("N86").Value = -8517.75

$ws.Range("I89").Value = 66672230
$ws.Range("J89").Value = 6271.75
$ws.Range("K89").Value = 333361150
$ws.Range("L89").Value = 31358.75
$ws.Range("M89").Value = -333355534
$ws.Range("N89").Value = -42590.75

$ws.Range("H94").Value = 1277.125
$ws.Range("I94").Value = 1277.125
$ws.Range("K94").Value = 1277.125
$ws.Range("M94").Value = -826.125

$ws.Range("H138").Value = 8489.35
$ws.Range("J138").Value = 4601.625
$ws.Range("L138").Value = 13804.875
$ws.Range("N138").Value = -24084.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 33333932
$ws.Range("I29").Value = 33333932
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 33333932
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -33333624
$ws.Range("N29").ClearContents()

$ws.Range("H113").Value = 40000
$ws.Range("J113").Value = 40000
$ws.Range("L113").Value = 40000
$ws.Range("N113").Value = -48678

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 5000
$ws.Range("I102").Value = 5000
$ws.Range("K102").Value = 5000
$ws.Range("M102").Value = -1755

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()

$ws.Range("H31").Value = 5147.095
$ws.Range("I31").Value = 3418.6
$ws.Range("J31").Value = 5687.25
$ws.Range("K31").Value = 3418.6
$ws.Range("L31").Value = 5687.25
$ws.Range("M31").Value = -3123.6
$ws.Range("N31").Value = -6277.25

$ws.Range("H34").Value = 5147.095
$ws.Range("I34").Value = 3418.6
$ws.Range("J34").Value = 5687.25
$ws.Range("K34").Value = 3418.6
$ws.Range("L34").Value = 5687.25
$ws.Range("M34").Value = -3216.6
$ws.Range("N34").Value = -6091.25

$ws.Range("H132").Value = 9011517
$ws.Range("I132").Value = 2381.8462
$ws.Range("K132").Value = 7145.5386
$ws.Range("M132").Value = -4615.5386

$ws.Range("H134").Value = 1629.4073
$ws.Range("I134").Value = 1695.6086
$ws.Range("J134").Value = 1248.75
$ws.Range("K134").Value = 5086.825800000001
$ws.Range("L134").Value = 3746.25
$ws.Range("M134").Value = -2551.825800000001
$ws.Range("N134").Value = -8816.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 367
$ws.Range("I31").Value = 367
$ws.Range("K31").Value = 1101
$ws.Range("M31").Value = -813

$ws.Range("H34").Value = 2949.0908
$ws.Range("J34").Value = 14000
$ws.Range("L34").Value = 42000
$ws.Range("N34").Value = -42168

$ws.Range("H68").Value = 2408.2778
$ws.Range("J68").Value = 2799.8572
$ws.Range("L68").Value = 8399.571599999999
$ws.Range("N68").Value = -10021.5716

$ws.Range("H71").Value = 2408.2778
$ws.Range("J71").Value = 2799.8572
$ws.Range("L71").Value = 25198.7148
$ws.Range("N71").Value = -33310.7148

$ws.Range("H113").Value = 695.1111
$ws.Range("J113").Value = 711.63635
$ws.Range("L113").Value = 2134.90905
$ws.Range("N113").Value = -6474.90905

$ws.Range("H138").Value = 3208.1738
$ws.Range("I138").Value = 3212.5264
$ws.Range("K138").Value = 9637.5792
$ws.Range("M138").Value = -4497.5792

$ws.Range("H141").Value = 2975.8333
$ws.Range("J141").Value = 2937.5
$ws.Range("L141").Value = 8812.5
$ws.Range("N141").Value = -19172.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 572357.1
$ws.Range("I3").Value = 700
$ws.Range("K3").Value = 700
$ws.Range("M3").Value = -584

$ws.Range("H70").Value = 4501.76
$ws.Range("I70").Value = 4251
$ws.Range("J70").Value = 5295.8335
$ws.Range("K70").Value = 4251
$ws.Range("L70").Value = 5295.8335
$ws.Range("M70").Value = -3981
$ws.Range("N70").Value = -5835.8335

$ws.Range("H73").Value = 4501.76
$ws.Range("I73").Value = 4251
$ws.Range("J73").Value = 5295.8335
$ws.Range("K73").Value = 4251
$ws.Range("L73").Value = 5295.8335
$ws.Range("M73").Value = -3315
$ws.Range("N73").Value = -7167.8335

$ws.Range("H80").Value = 74408.836
$ws.Range("I80").Value = 30035.191
$ws.Range("K80").Value = 30035.191
$ws.Range("M80").Value = -29037.191

$ws.Range("H83").Value = 74408.836
$ws.Range("I83").Value = 30035.191
$ws.Range("K83").Value = 150175.955
$ws.Range("M83").Value = -145183.955

$ws.Range("H113").Value = 2314.9395
$ws.Range("I113").Value = 2021.3448
$ws.Range("K113").Value = 2021.3448
$ws.Range("M113").Value = 148.6551999999999

$ws.Range("H132").Value = 10606809
$ws.Range("I132").Value = 5527.7144
$ws.Range("K132").Value = 16583.1432
$ws.Range("M132").Value = -14053.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7105.4443
$ws.Range("I40").Value = 8259.799999999999
$ws.Range("K40").Value = 8259.799999999999
$ws.Range("M40").Value = -8123.799999999999

$ws.Range("H68").Value = 3031.0322
$ws.Range("J68").Value = 3643.6667
$ws.Range("L68").Value = 3643.6667
$ws.Range("N68").Value = -5141.6667

$ws.Range("H71").Value = 3031.0322
$ws.Range("J71").Value = 3643.6667
$ws.Range("L71").Value = 18218.3335
$ws.Range("N71").Value = -25706.3335

$ws.Range("H87").Value = 500037500
$ws.Range("J87").Value = 500037500
$ws.Range("L87").Value = 500037500
$ws.Range("N87").Value = -500039746

$ws.Range("H90").Value = 500037500
$ws.Range("J90").Value = 500037500
$ws.Range("L90").Value = 1500112500
$ws.Range("N90").Value = -1500123732

$ws.Range("H122").Value = 3147.516
$ws.Range("I122").Value = 2903.6086
$ws.Range("J122").Value = 3848.75
$ws.Range("K122").Value = 8710.825800000001
$ws.Range("L122").Value = 11546.25
$ws.Range("M122").Value = -6260.825800000001
$ws.Range("N122").Value = -16446.25

$ws.Range("H136").Value = 23603.666
$ws.Range("I136").Value = 25932
$ws.Range("K136").Value = 77796
$ws.Range("M136").Value = -75246

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 9999.5
$ws.Range("I32").Value = 9999.5
$ws.Range("K32").Value = 9999.5
$ws.Range("M32").Value = -9682.5

$ws.Range("H100").Value = 1140.2
$ws.Range("I100").Value = 1103.3846
$ws.Range("K100").Value = 2206.7692
$ws.Range("M100").Value = -1665.7692

$ws.Range("H122").Value = 2269.7666
$ws.Range("I122").Value = 1775.7037
$ws.Range("K122").Value = 5327.1111
$ws.Range("M122").Value = -2877.1111

$ws.Range("H125").Value = 76993.25
$ws.Range("J125").Value = 76993.25
$ws.Range("L125").Value = 76993.25
$ws.Range("N125").Value = -86833.25

$ws.Range("H132").Value = 1321.341
$ws.Range("I132").Value = 1046.625
$ws.Range("J132").Value = 4068.5
$ws.Range("K132").Value = 3139.875
$ws.Range("L132").Value = 12205.5
$ws.Range("M132").Value = -609.875
$ws.Range("N132").Value = -17265.5
